$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "<there>"
$ws.Range("C2").Value = 28

# Row 3
$ws.Range("C3").Value = 36

# Row 4
$ws.Range("B4").Value = "<no>"
$ws.Range("C4").Value = 40

# Row 5
$ws.Range("C5").Value = 28

# Row 6
$ws.Range("B6").Value = "<is>"
$ws.Range("C6").Value = 29

# Row 7
$ws.Range("C7").Value = 37

# Row 8
$ws.Range("B8").Value = "<nimber>"
$ws.Range("C8").Value = 37

# Row 9
$ws.Range("B9").Value = "<coop>"
$ws.Range("C9").Value = 35

# Row 10
$ws.Range("B10").Value = "<canon>"
$ws.Range("C10").Value = 28

# Row 11
$ws.Range("B11").Value = "<eight>"
$ws.Range("C11").Value = 36

# Row 12
$ws.Range("C12").Value = 35

# Row 13
$ws.Range("C13").Value = 31

# Row 15
$ws.Range("C15").Value = 34

# Row 16
$ws.Range("B16").Value = "<nremo>"
$ws.Range("C16").Value = 27

# Row 17
$ws.Range("B17").Value = "<escat>"
$ws.Range("C17").Value = 38

# Row 18
$ws.Range("B18").Value = "<what>"
$ws.Range("C18").Value = 26
